$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header for column O
$ws.Range("O1").Value = "fy_str_unA"

# New column O: quarterly annualization fix = K / 4, for data rows 2 through 97.
# Filled in stages matching the existing L/M/N shared-formula boundaries
# (row 2 stand-alone, then 3:66, then 67:97) so the groups line up the same way.
$ws.Range("O2").Formula = "=K2/4"
$ws.Range("O3:O66").Formula = "=K3/4"
$ws.Range("O67:O97").Formula = "=K67/4"

# Select O2 like in the target sheet view
$ws.Range("O2").Select()
